$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new weekly record as row 28, shifting existing rows 28-67
# down to 29-68 (dimension grows from A1:R67 to A1:R68).
$ws.Rows.Item(28).Insert()

$ws.Range("A28").Value = 9
$ws.Range("B28").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C28").Value = "Metropolitana"
$ws.Range("D28").Value = 44483
$ws.Range("E28").Value = 13
$ws.Range("F28").Value = 100112022
$ws.Range("G28").Value = "Arveja Verde"
$ws.Range("H28").Value = "Perfection"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 34
$ws.Range("K28").Value = 24000
$ws.Range("L28").Value = 25000
$ws.Range("M28").Value = 24500
$ws.Range("N28").Value = "$/malla 25 kilos"
$ws.Range("O28").Value = "Provincia de Huasco"
$ws.Range("P28").Value = 980
$ws.Range("Q28").Value = 25
$ws.Range("R28").Value = "Hortaliza"
